$wb = $excel.ActiveWorkbook

# Sheet "block_0"
$ws1 = $wb.Worksheets.Item("block_0")
$ws1.Range("C2").Value = 0.2412331314783636
$ws1.Range("D2").Value = 522.7384330412233
$ws1.Range("E2").Value = 0.1567800643885285
$ws1.Range("F2").Value = 187.2378100545062
$ws1.Rows.Item(3).Delete()

# Sheet "block_1"
$ws2 = $wb.Worksheets.Item("block_1")
$ws2.Range("C2").Value = 0.0112197145813409
$ws2.Range("D2").Value = 1524.443749419503
$ws2.Range("E2").Value = 0.005926851996257251
$ws2.Range("F2").Value = 469.8066630824756
$ws2.Rows.Item(3).Delete()
